# Updated cryptos list on Tue Sep  3 04:31:05 UTC 2024 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns with the latest scrape;
# rows 30/31 (PancakeSwap / Monero) also swap ranking positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '58.845.29'
$ws.Range('E2').Value = '  +2.21%  '
# Row 3
$ws.Range('D3').Value = '2.490.02'
$ws.Range('E3').Value = '  +2.45%  '
# Row 4
$ws.Range('E4').Value = '  +0.13%  '
# Row 5
$ws.Range('D5').Value = '''533.04'
$ws.Range('E5').Value = '  +3.96%  '
# Row 6
$ws.Range('D6').Value = '''134.88'
$ws.Range('E6').Value = '  +4.33%  '
# Row 7
$ws.Range('E7').Value = '  +0.30%  '
# Row 8
$ws.Range('D8').Value = '''0.565'
$ws.Range('E8').Value = '  +3.02%  '
# Row 9
$ws.Range('D9').Value = '2.512.90'
$ws.Range('E9').Value = '  +2.92%  '
# Row 10
$ws.Range('D10').Value = '''0.0994'
$ws.Range('E10').Value = '  +4.79%  '
# Row 11
$ws.Range('E11').Value = '  -1.44%  '
# Row 12
$ws.Range('D12').Value = '''5.25'
$ws.Range('E12').Value = '  +1.56%  '
# Row 13
$ws.Range('D13').Value = '''0.334'
$ws.Range('E13').Value = '  +1.21%  '
# Row 14
$ws.Range('D14').Value = '2.939.89'
$ws.Range('E14').Value = '  +2.64%  '
# Row 15
$ws.Range('D15').Value = '58.920.96'
$ws.Range('E15').Value = '  +2.49%  '
# Row 16
$ws.Range('D16').Value = '''22.39'
$ws.Range('E16').Value = '  +2.95%  '
# Row 17
$ws.Range('E17').Value = '  +3.21%  '
# Row 18
$ws.Range('D18').Value = '2.515.91'
$ws.Range('E18').Value = '  +3.08%  '
# Row 19
$ws.Range('D19').Value = '''10.66'
$ws.Range('E19').Value = '  +2.32%  '
# Row 20
$ws.Range('D20').Value = '''4.23'
$ws.Range('E20').Value = '  +3.40%  '
# Row 21
$ws.Range('D21').Value = '''321.46'
$ws.Range('E21').Value = '  +2.03%  '
# Row 22
$ws.Range('D22').Value = '''6.12'
$ws.Range('E22').Value = '  +8.66%  '
# Row 23
$ws.Range('D23').Value = '''0.995'
$ws.Range('E23').Value = '  -0.46%  '
# Row 24
$ws.Range('D24').Value = '''65.84'
$ws.Range('E24').Value = '  +3.79%  '
# Row 25
$ws.Range('E25').Value = '  +1.12%  '
# Row 26
$ws.Range('D26').Value = '''0.997'
$ws.Range('E26').Value = '  +0.14%  '
# Row 27
$ws.Range('E27').Value = '  +1.11%  '
# Row 28
$ws.Range('D28').Value = '''7.47'
$ws.Range('E28').Value = '  +3.71%  '
# Row 29
$ws.Range('D29').Value = '0.0₃0764'
$ws.Range('E29').Value = '  +5.96%  '
# Row 30
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '''171.47'
$ws.Range('E30').Value = '  +0.71%  '
# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''1.75'
$ws.Range('E31').Value = '  +4.72%  '
# Row 32
$ws.Range('E32').Value = '  +3.79%  '
# Row 33
$ws.Range('D33').Value = '''6.29'
$ws.Range('E33').Value = '  +1.14%  '
# Row 34
$ws.Range('D34').Value = '''0.999'
$ws.Range('E34').Value = '  +0.00%  '
# Row 36
$ws.Range('D36').Value = '''18.18'
$ws.Range('E36').Value = '  +2.92%  '
# Row 37
$ws.Range('D37').Value = '''1.25'
$ws.Range('E37').Value = '  -1.64%  '
# Row 38
$ws.Range('D38').Value = '''3.97'
$ws.Range('E38').Value = '  +1.65%  '
# Row 39
$ws.Range('E39').Value = '  +4.32%  '
# Row 40
$ws.Range('D40').Value = '''36.70'
$ws.Range('E40').Value = '  +1.22%  '
# Row 41
$ws.Range('D41').Value = '''0.783'
$ws.Range('E41').Value = '  +1.24%  '
# Row 42
$ws.Range('D42').Value = '''279.48'
$ws.Range('E42').Value = '  +2.71%  '
# Row 43
$ws.Range('D43').Value = '''3.48'
$ws.Range('E43').Value = '  +3.42%  '
# Row 44
$ws.Range('D44').Value = '''5.07'
$ws.Range('E44').Value = '  +3.88%  '
# Row 45
$ws.Range('D45').Value = '''131.80'
$ws.Range('E45').Value = '  +9.54%  '
# Row 46
$ws.Range('D46').Value = '''0.593'
$ws.Range('E46').Value = '  +1.66%  '
# Row 47
$ws.Range('D47').Value = '''0.0929'
$ws.Range('E47').Value = '  +2.47%  '
# Row 48
$ws.Range('E48').Value = '  +5.56%  '
# Row 49
$ws.Range('E49').Value = '  +4.01%  '
# Row 50
$ws.Range('D50').Value = '''17.06'
$ws.Range('E50').Value = '  +3.52%  '
# Row 51
$ws.Range('D51').Value = '1.758.61'
$ws.Range('E51').Value = '  +3.14%  '
